$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: remove the "License Information" Heading2 paragraph entirely.
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "License Information") {
        $p.Range.Delete()
        $found = $true
        break
    }
}

# ---------------------------------------------------------------------------
# Step 2: remove the "This PDF version is provided under the same license."
# paragraph entirely (its content is dropped, not merged anywhere).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "This PDF version is provided under the same license.") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Step 3: rewrite the license/credits paragraph (the one that used to start
# with the bold "Questions de Traduction (unfoldingWord)" run).
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.StartsWith(" Questions de Traduction (unfoldingWord)") -or $t.StartsWith("Questions de Traduction (unfoldingWord) (French) is based on")) {
        $target = $p
        break
    }
}

$r = $target.Range
$start = $r.Start
$end = $r.End - 1
$clear = $d.Range($start, $end)
$clear.Text = ""

$pos = $start

$segments = @(
    @{ Text = "unfoldingWord® Translation Questions"; Bold = $true },
    @{ Text = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. "; Bold = $false },
    @{ Text = "unfoldingWord® Translation Questions"; Bold = $false },
    @{ Text = " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from "; Bold = $false },
    @{ Text = "unfoldingWord® Translation Questions"; Bold = $false },
    @{ Text = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"; Bold = $false }
)

foreach ($seg in $segments) {
    $segStart = $pos
    $ins = $d.Range($pos, $pos)
    $ins.InsertAfter($seg.Text)
    $segEnd = $segStart + $seg.Text.Length
    $fr = $d.Range($segStart, $segEnd)
    if ($seg.Bold) {
        $fr.Font.Bold = 1
    } else {
        $fr.Font.Bold = 0
    }
    $pos = $segEnd
}

Write-Output "Final paragraph text: $($target.Range.Text)"
Write-Output "Paragraph count: $($d.Paragraphs.Count)"
